$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# squats row: drop the one-off "head(squats,1)" formula text in favour of "sticky"
$ws.Range("J6").Value = "sticky"

# Insert a new survey question (Jumping jacks) right after squats
$ws.Rows.Item(7).Insert()
$ws.Range("F7").Value = "Jumping jacks"
$ws.Range("D7").Value = "jumping_jacks"
$ws.Range("C7").Value = "number 0,1000,1"
$ws.Range("J7").Value = "sticky"
$ws.Rows.Item(7).RowHeight = 30

# Remove the old "submit" and trailing "note" rows (Go on! / Good work, chap!)
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(14).Delete()

# Add a new file-upload question for a progress photo
$ws.Range("C15").Value = "file"
$ws.Range("F15").Value = "Photograph your abs!"
$ws.Range("D15").Value = "abs_image"

# Update the selected cell to reflect where editing left off
$ws.Range("D15").Select()
